# Auto-generated edit script reproducing the cryptos.xlsx data refresh
# (GitHub Actions commit "Updated cryptos list on Sat Aug  3 13:46:45 UTC 2024").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.747.74"
$ws.Range("E2").Value = "  -5.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.988.35"
$ws.Range("E3").Value = "  -5.36%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.72"
$ws.Range("E5").Value = "  -5.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.47"
$ws.Range("E6").Value = "  -9.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.004.19"
$ws.Range("E9").Value = "  -5.32%  "
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -7.70%  "
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.510.65"
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.801.79"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.98"
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.996.87"
$ws.Range("E17").Value = "  -4.85%  "
$ws.Range("E18").Value = "  -6.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.16"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.06"
$ws.Range("E20").Value = "  -4.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.05"
$ws.Range("E21").Value = "  -9.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").Value = "  -5.71%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.11"
$ws.Range("E24").Value = "  -4.37%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.115.17"
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.469"
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0932"
$ws.Range("E29").Value = "  -11.41%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.22"
$ws.Range("E30").Value = "  -11.01%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.72"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.43"
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.99"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.98"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.59"
$ws.Range("E36").Value = "  -7.44%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.07"
$ws.Range("E37").Value = "  -6.54%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.28"
$ws.Range("E38").Value = "  -7.22%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  -8.84%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.67"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.415.26"
$ws.Range("E41").Value = "  -7.85%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").Value = "  -6.83%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.07"
$ws.Range("E43").Value = "  -8.42%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.673"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0591"
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.19"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0245"
$ws.Range("E48").Value = "  -5.71%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0953"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.67"
$ws.Range("E50").Value = "  -8.52%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "266.55"
$ws.Range("E51").Value = "  -9.37%  "
